$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 997
$ws.Range("B2").Value = 995
$ws.Range("C2").Value = 995
$ws.Range("D2").Value = 995
$ws.Range("E2").Value = 997
$ws.Range("F2").Value = 999
$ws.Range("G2").Value = 996
$ws.Range("H2").Value = 998
